# The template's "body-center / body-bold-center" paragraph originally held
# two template blocks back-to-back on one line:
#   {#is_bodycenter}{message}{/is_bodycenter}{#is_bodyboldcenter}{message}{/is_ bodyboldcenter}
# The edit splits that single paragraph into two paragraphs, breaking right
# between "{/is_bodycenter}" and "{#is_bodyboldcenter}" (i.e. the user put
# their cursor there and pressed Enter). The new paragraph inherits the same
# paragraph formatting (centered) as the original.

$d = $word.ActiveDocument

$marker = "{/is_bodycenter}"

# Locate the exact boundary between the two template blocks without
# replacing any text (a Find/Replace spanning this point would also delete
# the bookmark that sits here), then collapse to a single point right after
# "{/is_bodycenter}" and insert a paragraph break there.
$searchRange = $d.Content
$found = $searchRange.Find.Execute($marker + "{#is_bodyboldcenter}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $splitPoint = $searchRange.Start + $marker.Length
    $insertionPoint = $d.Range($splitPoint, $splitPoint)
    $insertionPoint.InsertParagraphAfter()
}
